# Swahili -> English translation edits for
# "05_Prisoners and candies - subtitles (corrected).docx"
#
# We locate each piece of text with Find.Execute and then assign the
# found Range's .Text property directly (instead of passing the
# replacement through Find.Execute's Replace argument) so that Word's
# "replace straight quotes with smart quotes" AutoCorrect/AutoFormat
# behaviour does not mangle the straight apostrophes used throughout
# the target English text.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
    $rng.Text = $replace
}

Replace-Text "Wafungwa na peremende - manukuu:" "Prisoners and candies - subtitles:"

Replace-Text "**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino" "**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino"

Replace-Text "[Muziki]" "[Music]"

Replace-Text "wanahisabati wanne mkali wanachukuliwa" "four bright mathematicians are taken into"

Replace-Text "chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu" "custody and put in jail because they tried"

Replace-Text "kumshawishi mwanamke mzee kuwa Goedel's" "to convince an old lady that the Goedel's"

Replace-Text "nadharia za kutokamilika ni kweli. Kila" "incompleteness theorems are true. Every"

Replace-Text "mtaalamu wa hisabati ana kiini chake ambacho sisi" "mathematician has his own cell that we"

Replace-Text "inaweza kuhesabu na nambari kutoka 1 hadi 4." "can enumerate with a number from 1 to 4."

Replace-Text "kabla ya kuingia kwenye seli fulani" "before entering the cell a certain"

Replace-Text "idadi ya peremende kubwa kuliko " "number of candies greater than "

Replace-Text "e: AU " "e: OR "

Replace-Text "SAWA NA" "EQUAL TO"

# Remove the standalone single-space run between ")" and " 1 ni", collapsing
# the double space into a single space (that run disappears entirely in the
# edited document).
Replace-Text ") " ")"

Replace-Text " 1 ni" " 1 is"

Replace-Text "wanapewa kila mtaalamu wa hisabati na wao" "given to every mathematician and they"

Replace-Text "wanaambiwa wana peremende 11 kwa jumla." "are told they have 11 candies in total."

Replace-Text "lakini kila mtu anajua idadi yake tu" "but everyone knows only his number of"

Replace-Text "pipi na jumla. 1 na sio" "candies and the total. 1 and is not"

Replace-Text "kuruhusiwa kuuliza nambari zingine." "allowed to ask for the others number."

Replace-Text "kisha mwanahisabati wa kwanza anauliza" "then the first mathematician asks the"

Replace-Text "pili: 'namba 2 unajua kama wewe" "second: 'number 2 do you know if you"

Replace-Text "kuwa na peremende nyingi kuliko mimi?' ya pili" "have more candies than me?' the second"

Replace-Text "mwanahisabati anajibu hana. Kisha yeye" "mathematician answers he doesn't. Then he"

Replace-Text "anauliza kwa nambari 3: 'unajua kama unayo" "asks to number 3: 'do you know if you have"

Replace-Text "pipi zaidi kuliko mimi?'" "more candy than me?'"

Replace-Text "mwanahisabati wa tatu anajibu: 'hapana niko" "the third mathematician answers: 'no I'm"

Replace-Text "samahani sifanyi'. Katika hatua hii ya nne" "sorry I don't'. At this point the fourth"

Replace-Text "mtaalamu wa hisabati anasema: 'jamani mnafahamu" "mathematician says: 'hey guys you know"

Replace-Text "nini, najua hasa pipi ngapi" "what, I know exactly how many candies"

Replace-Text "kila mtu ana hapa'. Cha kushangaza hata" "everyone has here'. Surprisingly even the"

Replace-Text "wanahisabati wengine watatu wanasema hivyo sasa" "other three mathematicians say that now"

Replace-Text "wanajua kila mtu ana pipi ngapi" "they know how many candies everyone has"

Replace-Text "kwa hivyo swali ni: unaweza kujua" "so the question is: can you figure out"

Replace-Text "idadi ya pipi kila mfungwa ana" "the number of candies every prisoner has"

# Second occurrence of the music cue, at the end of the subtitles.
Replace-Text "[Muziki]" "[Music]"
